$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 41

# Copy the formatting (alignment etc.) of the previous data row onto the new
# row first, so A41:C41 end up using the same style as the rest of the table.
$ws.Range("A40:C40").Copy()
$ws.Range("A41:C41").PasteSpecial(-4122)

# Column A stores "dates" as plain text (matching every other row, which is
# inlineStr/text, not a real date). Mark the cell as Text *before* writing
# the value so the engine doesn't auto-convert "2025/12/20" into a date
# serial number.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025/12/20"
$ws.Cells.Item($row, 2).Value = "逃离鸭科夫"
$ws.Cells.Item($row, 3).Value = 1357
